$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Monthly Budget"

$range = $ws.Range("B13:B16")
$range.Font.Bold = $true
$range.Font.Size = 12
$range.Font.ThemeColor = 1
$range.Interior.Color = 5287936
$range.HorizontalAlignment = -4152
$ws.Rows("13:16").RowHeight = 16
